# Apply the price / 24h-volume-change refresh captured by the scheduled
# GitHub Actions "Updated cryptos list" commit.
#
# Column D ("Price") cells are stored as literal text in the workbook (many
# values use "." as a thousands separator, e.g. "29.391.98", and some need
# to keep insignificant trailing zeros, e.g. "0.01860"). Assigning a plain
# numeric-looking string to Range.Value lets Excel auto-coerce it to a
# number (losing the trailing zero / switching to scientific notation for
# very small values), so for the subset of new prices that would otherwise
# be parsed as numbers we force the cell to Text format first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textPriceCells = @(
    "D4", "D5", "D7", "D8", "D9", "D10", "D11", "D13", "D14", "D15",
    "D16", "D18", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27",
    "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D39",
    "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D49", "D50"
)
foreach ($addr in $textPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Cell-by-cell updates, in the same order the source table lists them.
$ws.Range("D2").Value = "29.391.98"
$ws.Range("E2").Value = "  +1.80%  "
$ws.Range("D3").Value = "1.851.80"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "245.05"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "0.3055"
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").Value = "0.07642"
$ws.Range("E9").Value = "  -0.81%  "
$ws.Range("D10").Value = "23.41"
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("D11").Value = "0.07728"
$ws.Range("E11").Value = "  -1.02%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.851.02"
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "5.134"
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("D14").Value = "0.6917"
$ws.Range("E14").Value = "  +1.38%  "
$ws.Range("D15").Value = "90.75"
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("D16").Value = "6.302"
$ws.Range("E16").Value = "  -1.86%  "
$ws.Range("D17").Value = "29.398.46"
$ws.Range("E17").Value = "  +1.86%  "
$ws.Range("D18").Value = "0.000008254"
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("D19").Value = "2.095.39"
$ws.Range("E19").Value = "  +0.99%  "
$ws.Range("D20").Value = "235.86"
$ws.Range("E20").Value = "  -2.73%  "
$ws.Range("D21").Value = "12.69"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").Value = "0.9995"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "7.636"
$ws.Range("E23").Value = "  +2.49%  "
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").Value = "0.1472"
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("D26").Value = "8.951"
$ws.Range("E26").Value = "  +1.69%  "
$ws.Range("D27").Value = "160.18"
$ws.Range("E27").Value = "  +1.41%  "
$ws.Range("D28").Value = "18.17"
$ws.Range("E28").Value = "  -0.44%  "
$ws.Range("D29").Value = "1.528"
$ws.Range("E29").Value = "  -1.03%  "
$ws.Range("D30").Value = "4.249"
$ws.Range("E30").Value = "  +0.64%  "
$ws.Range("D31").Value = "4.131"
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("D32").Value = "1.199"
$ws.Range("E32").Value = "  +0.31%  "
$ws.Range("D33").Value = "0.05221"
$ws.Range("E33").Value = "  +2.51%  "
$ws.Range("D34").Value = "0.7711"
$ws.Range("E34").Value = "  -0.53%  "
$ws.Range("D35").Value = "1.871"
$ws.Range("E35").Value = "  +1.04%  "
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("D37").Value = "2.685"
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("D38").Value = "1.315.36"
$ws.Range("E38").Value = "  +7.20%  "
$ws.Range("D39").Value = "0.01860"
$ws.Range("E39").Value = "  +0.61%  "
$ws.Range("E40").Value = "  +0.91%  "
$ws.Range("D41").Value = "0.9414"
$ws.Range("E41").Value = "  -1.28%  "
$ws.Range("D42").Value = "105.81"
$ws.Range("E42").Value = "  -3.00%  "
$ws.Range("D43").Value = "5.779"
$ws.Range("E43").Value = "  -1.71%  "
$ws.Range("D44").Value = "0.9998"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "9.741"
$ws.Range("E45").Value = "  +1.27%  "
$ws.Range("D46").Value = "1.997.15"
$ws.Range("E46").Value = "  +1.08%  "
$ws.Range("D47").Value = "0.5221"
$ws.Range("E47").Value = "  +1.25%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.00000000122"
$ws.Range("E48").Value = "  -0.61%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "1.774"
$ws.Range("E49").Value = "  +1.41%  "
$ws.Range("D50").Value = "62.81"
$ws.Range("E50").Value = "  -2.17%  "
